$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item("Sheet1")

# --- Sheet1: view changes (no longer the tab shown; scrolled/zoomed) ---
$sheet1.Activate()
$sheet1.Application.ActiveWindow.Zoom = 139
$sheet1.Range("A8").Select()
$sheet1.Application.ActiveWindow.ScrollRow = 8
$sheet1.Range("A37").Select()

# --- Add Sheet2 right after Sheet1 ---
$ws2 = $wb.Worksheets.Add($null, $sheet1)
$ws2.Name = "Sheet2"

# Header row
$ws2.Range("A2").Value = "config file"
$ws2.Range("B2").Value = "result"
$ws2.Range("C2").Value = "tensorboard"

# Row 3
$ws2.Range("A3").Value = "cifar10_20_dir_p3_se_conv4_kaiming_init"
$ws2.Range("A3").Font.Name = "JetBrains Mono"
$ws2.Range("A3").Font.Family = 3
$ws2.Range("A3").Font.Size = 16
$ws2.Range("B3").Value = "up to 70% then down to random"
$ws2.Range("C3").Value = "model_Conv4_n_cli_20_ds_split_dirichlet_ds_alpha_0.3_align_se_delta_None_init_type_kaiming_normal_same_init_True"
$ws2.Range("C3").Font.Name = "Calibri Light"
$ws2.Range("C3").Font.ThemeFont = 2
$ws2.Range("C3").Font.Size = 12
$ws2.Rows("3:3").RowHeight = 20

# Row 4
$ws2.Range("A4").Value = "cifar10_20_dir_p3_ae_conv4_default_init"
$ws2.Range("A4").Font.Name = "JetBrains Mono"
$ws2.Range("A4").Font.Family = 3
$ws2.Range("A4").Font.Size = 14
$ws2.Range("B4").Value = "85% at 324"
$ws2.Range("C4").Value = "model_Conv4_n_cli_20_ds_split_dirichlet_ds_alpha_0.3_align_ae_delta_None_init_type_default_same_init_True"
$ws2.Range("C4").Font.Name = "Calibri Light"
$ws2.Range("C4").Font.ThemeFont = 2
$ws2.Range("C4").Font.Size = 12
$ws2.Rows("4:4").RowHeight = 18

# Row 5
$ws2.Range("A5").Value = "cifar10_20_dir_p3_ae_conv4_kaiming"
$ws2.Range("B5").Value = "85% at 273"
$ws2.Range("C5").Value = "model_Conv4_n_cli_20_ds_split_dirichlet_ds_alpha_0.3_align_ae_delta_None_init_type_kaiming_normal_same_init_True/Accuracy/test/max_client_test_accuracy"
$ws2.Range("C5").Font.Name = "Calibri Light"
$ws2.Range("C5").Font.ThemeFont = 2
$ws2.Range("C5").Font.Size = 13
$ws2.Rows("5:5").RowHeight = 17

# Row 32
$ws2.Range("A32").Value = "s"

# Column widths
$ws2.Columns("A:A").ColumnWidth = 47.6640625
$ws2.Columns("B:B").ColumnWidth = 28
$ws2.Columns("C:C").ColumnWidth = 135.83203125

# Sheet2 view state: active sheet/tab, selection, scroll position
$ws2.Activate()
$ws2.Range("A32").Select()
$wb.Application.ActiveWindow.WindowState = $wb.Application.ActiveWindow.WindowState
